$d = $word.ActiveDocument

function Rename-InlinePicture($headerFooter, $newName) {
    $inlineShape = $headerFooter.Range.InlineShapes.Item(1)
    # InlineShape has no writable .Name in the Word object model, so the
    # picture is temporarily promoted to a floating Shape (which does
    # expose .Name -> the drawing's docPr/@name), renamed, then converted
    # straight back to an inline picture (wp:inline is restored, not left
    # as wp:anchor).
    $floating = $inlineShape.ConvertToShape()
    try {
        $floating.Name = $newName
    } finally {
        [void]$floating.ConvertToInlineShape()
    }
}

# Default (primary) footer: Pearson logo (docPr id=1 / cNvPr id=0) -> image2.png becomes image1.png
$sec = $d.Sections.First
$ftrDefault = $sec.Footers.Item(1)
if ($ftrDefault.Exists) {
    Rename-InlinePicture $ftrDefault "image1.png"
}

# First-page footer: Pearson logo (docPr id=2 / cNvPr id=0) -> image2.png becomes image1.png
$sec = $d.Sections.First
$ftrFirst = $sec.Footers.Item(2)
if ($ftrFirst.Exists) {
    Rename-InlinePicture $ftrFirst "image1.png"
}

# First-page header: BTec logo (docPr id=3 / cNvPr id=0) -> image1.jpg becomes image2.jpg
$sec = $d.Sections.First
$hdrFirst = $sec.Headers.Item(2)
if ($hdrFirst.Exists) {
    Rename-InlinePicture $hdrFirst "image2.jpg"
}
